# Insert a new price-record row at row 434 (a new weekly "Ajo" quote for
# Feria Lagunitas de Puerto Montt). This pushes the existing rows 434-490
# down to 435-491 and grows the used range from A1:R490 to A1:R491.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(434).Insert()

$ws.Cells.Item(434, 1).Value  = 4
$ws.Cells.Item(434, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(434, 3).Value  = "Los Lagos"
$ws.Cells.Item(434, 4).Value  = 45127
$ws.Cells.Item(434, 5).Value  = 10
$ws.Cells.Item(434, 6).Value  = 100112003
$ws.Cells.Item(434, 7).Value  = "Ajo"
$ws.Cells.Item(434, 8).Value  = "Chino"
$ws.Cells.Item(434, 9).Value  = "Primera"
$ws.Cells.Item(434, 10).Value = 120
$ws.Cells.Item(434, 11).Value = 22000
$ws.Cells.Item(434, 12).Value = 22000
$ws.Cells.Item(434, 13).Value = 22000
$ws.Cells.Item(434, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(434, 15).Value = "China"
$ws.Cells.Item(434, 16).Value = 2200
$ws.Cells.Item(434, 17).Value = 10
$ws.Cells.Item(434, 18).Value = "Hortaliza"
